# Refresh Universalis market-price derived columns (H:N) across all job sheets,
# mirroring the scheduled market-data runner output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1168
$ws.Range("J17").Value = 1168
$ws.Range("L17").Value = 3504
$ws.Range("N17").Value = -3840

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6945207
$ws.Range("I43").Value = 501
$ws.Range("J43").Value = 7937308
$ws.Range("K43").Value = 501
$ws.Range("L43").Value = 7937308
$ws.Range("M43").Value = -432
$ws.Range("N43").Value = -7937446

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 203.41667
$ws.Range("I55").Value = 143.875
$ws.Range("J55").Value = 322.5
$ws.Range("K55").Value = 143.875
$ws.Range("L55").Value = 322.5
$ws.Range("M55").Value = 70.125
$ws.Range("N55").Value = -750.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 577.5833
$ws.Range("I96").Value = 564.6667
$ws.Range("J96").Value = 590.5
$ws.Range("K96").Value = 1694.0001
$ws.Range("L96").Value = 1771.5
$ws.Range("M96").Value = -321.0001
$ws.Range("N96").Value = -4517.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2441.9473
$ws.Range("I106").Value = 2441.9473
$ws.Range("K106").Value = 2441.9473
$ws.Range("M106").Value = -1810.9473

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 887
$ws.Range("I125").Value = 576
$ws.Range("J125").Value = 1042.5
$ws.Range("K125").Value = 5184
$ws.Range("L125").Value = 9382.5
$ws.Range("M125").Value = -2724
$ws.Range("N125").Value = -14302.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 8553853
$ws.Range("I132").Value = 13891986
$ws.Range("K132").Value = 41675958
$ws.Range("M132").Value = -41673428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 344.65
$ws.Range("I135").Value = 281.3846
$ws.Range("J135").Value = 462.14285
$ws.Range("K135").Value = 2532.4614
$ws.Range("L135").Value = 4159.28565
$ws.Range("M135").Value = 2.538600000000315
$ws.Range("N135").Value = -9229.28565

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 988
$ws.Range("I137").Value = 832.6667
$ws.Range("J137").Value = 1081.2
$ws.Range("K137").Value = 2498.0001
$ws.Range("L137").Value = 3243.6
$ws.Range("M137").Value = 51.9998999999998
$ws.Range("N137").Value = -8343.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 751485.25
$ws.Range("I138").Value = 1258.2
$ws.Range("J138").Value = 1013192.4
$ws.Range("K138").Value = 3774.6
$ws.Range("L138").Value = 3039577.2
$ws.Range("M138").Value = 1365.4
$ws.Range("N138").Value = -3049857.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 68326.664
$ws.Range("J139").Value = 68326.664
$ws.Range("L139").Value = 68326.664
$ws.Range("N139").Value = -78606.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1490.52
$ws.Range("I61").Value = 1403.15
$ws.Range("J61").Value = 1840
$ws.Range("K61").Value = 1403.15
$ws.Range("L61").Value = 1840
$ws.Range("M61").Value = -1191.15
$ws.Range("N61").Value = -2264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1113.5
$ws.Range("I74").Value = 1014.6
$ws.Range("J74").Value = 1360.75
$ws.Range("K74").Value = 1014.6
$ws.Range("L74").Value = 1360.75
$ws.Range("M74").Value = -140.6
$ws.Range("N74").Value = -3108.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1113.5
$ws.Range("I77").Value = 1014.6
$ws.Range("J77").Value = 1360.75
$ws.Range("K77").Value = 5073
$ws.Range("L77").Value = 6803.75
$ws.Range("M77").Value = -705
$ws.Range("N77").Value = -15539.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2479.9375
$ws.Range("I88").Value = 2051.8333
$ws.Range("K88").Value = 2051.8333
$ws.Range("M88").Value = -1645.8333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2479.9375
$ws.Range("I91").Value = 2051.8333
$ws.Range("K91").Value = 2051.8333
$ws.Range("M91").Value = -647.8332999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1674.1666
$ws.Range("I110").Value = 1388.1333
$ws.Range("J110").Value = 3104.3333
$ws.Range("K110").Value = 1388.1333
$ws.Range("L110").Value = 3104.3333
$ws.Range("M110").Value = 656.8667
$ws.Range("N110").Value = -7194.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1170.5294
$ws.Range("I122").Value = 1126.0769
$ws.Range("K122").Value = 3378.2307
$ws.Range("M122").Value = -928.2307000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1490.52
$ws.Range("I136").Value = 1403.15
$ws.Range("J136").Value = 1840
$ws.Range("K136").Value = 4209.450000000001
$ws.Range("L136").Value = 5520
$ws.Range("M136").Value = -1659.450000000001
$ws.Range("N136").Value = -10620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1866.9546
$ws.Range("I20").Value = 1850.4615
$ws.Range("J20").Value = 1890.7778
$ws.Range("K20").Value = 1850.4615
$ws.Range("L20").Value = 1890.7778
$ws.Range("M20").Value = -1603.4615
$ws.Range("N20").Value = -2384.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3874.606
$ws.Range("I86").Value = 4149.75
$ws.Range("K86").Value = 4149.75
$ws.Range("M86").Value = -3026.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3874.606
$ws.Range("I89").Value = 4149.75
$ws.Range("K89").Value = 20748.75
$ws.Range("M89").Value = -15132.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1528.2778
$ws.Range("I31").Value = 1521.1177
$ws.Range("J31").Value = 1650
$ws.Range("K31").Value = 1521.1177
$ws.Range("L31").Value = 1650
$ws.Range("M31").Value = -1226.1177
$ws.Range("N31").Value = -2240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1528.2778
$ws.Range("I34").Value = 1521.1177
$ws.Range("J34").Value = 1650
$ws.Range("K34").Value = 1521.1177
$ws.Range("L34").Value = 1650
$ws.Range("M34").Value = -1319.1177
$ws.Range("N34").Value = -2054

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5885326
$ws.Range("I62").Value = 3002.625
$ws.Range("K62").Value = 3002.625
$ws.Range("M62").Value = -2378.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5885326
$ws.Range("I65").Value = 3002.625
$ws.Range("K65").Value = 15013.125
$ws.Range("M65").Value = -11893.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5788
$ws.Range("J80").Value = 5788
$ws.Range("L80").Value = 17364
$ws.Range("N80").Value = -19236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5788
$ws.Range("J83").Value = 5788
$ws.Range("L83").Value = 52092
$ws.Range("N83").Value = -61452

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13640935
$ws.Range("I70").Value = 11908836
$ws.Range("J70").Value = 16672109
$ws.Range("K70").Value = 11908836
$ws.Range("L70").Value = 16672109
$ws.Range("M70").Value = -11908566
$ws.Range("N70").Value = -16672649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13640935
$ws.Range("I73").Value = 11908836
$ws.Range("J73").Value = 16672109
$ws.Range("K73").Value = 11908836
$ws.Range("L73").Value = 16672109
$ws.Range("M73").Value = -11907900
$ws.Range("N73").Value = -16673981

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4222.4
$ws.Range("I102").Value = 4678
$ws.Range("K102").Value = 4678
$ws.Range("M102").Value = -3056

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 992.4211
$ws.Range("I16").Value = 1004.1875
$ws.Range("K16").Value = 1004.1875
$ws.Range("M16").Value = -834.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2100
$ws.Range("J100").Value = 2333.3333
$ws.Range("L100").Value = 2333.3333
$ws.Range("N100").Value = -3415.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 129173.875
$ws.Range("I132").Value = 4629.6665
$ws.Range("J132").Value = 203900.4
$ws.Range("K132").Value = 13888.9995
$ws.Range("L132").Value = 611701.2
$ws.Range("M132").Value = -11358.9995
$ws.Range("N132").Value = -616761.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 46208.668
$ws.Range("J133").Value = 46208.668
$ws.Range("L133").Value = 46208.668
$ws.Range("N133").Value = -51268.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9363.23
$ws.Range("I136").Value = 12513
$ws.Range("J136").Value = 2276.25
$ws.Range("K136").Value = 37539
$ws.Range("L136").Value = 6828.75
$ws.Range("M136").Value = -34989
$ws.Range("N136").Value = -11928.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 7000
$ws.Range("J108").Value = 7000
$ws.Range("L108").Value = 7000
$ws.Range("N108").Value = -14680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 391.67856
$ws.Range("I113").Value = 348.83334
$ws.Range("J113").Value = 468.8
$ws.Range("K113").Value = 1046.50002
$ws.Range("L113").Value = 1406.4
$ws.Range("M113").Value = 1123.49998
$ws.Range("N113").Value = -5746.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 703
$ws.Range("I136").Value = 304.33334
$ws.Range("K136").Value = 913.0000200000001
$ws.Range("M136").Value = 1636.99998
